$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").Value = 25
$ws.Range("B11").Value = "pikachu"
$ws.Range("C11").Value = "electric"
